# Updated cryptocurrency price / 1h-volume data (and the Stellar/Kaspa row
# swap at rows 36-37) per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells are plain, unstyled text cells ("General" format, no
# cell style) in the source workbook. Assigning a numeric-looking string
# straight to .Value lets Excel auto-convert it to a real number (e.g.
# "13.60" -> 13.6, dropping the trailing zero), which would not match the
# source data. Forcing text format first preserves the literal string; since
# none of these cells carry any other formatting, ClearFormats() afterwards
# safely restores the original (unstyled) appearance.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "43.981.48"
Set-TextValue $ws.Range("E2") "  +0.07%  "

Set-TextValue $ws.Range("D3") "2.237.87"
Set-TextValue $ws.Range("E3") "  -0.21%  "

Set-TextValue $ws.Range("E4") "  +0.14%  "

Set-TextValue $ws.Range("D5") "306.59"
Set-TextValue $ws.Range("E5") "  -3.91%  "

Set-TextValue $ws.Range("D6") "95.02"
Set-TextValue $ws.Range("E6") "  -5.02%  "

Set-TextValue $ws.Range("E7") "  -0.62%  "

Set-TextValue $ws.Range("E8") "  +0.20%  "

Set-TextValue $ws.Range("D9") "0.521"
Set-TextValue $ws.Range("E9") "  -4.06%  "

Set-TextValue $ws.Range("D10") "34.84"
Set-TextValue $ws.Range("E10") "  -4.77%  "

Set-TextValue $ws.Range("D11") "0.0808"
Set-TextValue $ws.Range("E11") "  -1.95%  "

Set-TextValue $ws.Range("D12") "7.22"
Set-TextValue $ws.Range("E12") "  -3.50%  "

Set-TextValue $ws.Range("E13") "  -1.01%  "

Set-TextValue $ws.Range("D14") "2.579.83"
Set-TextValue $ws.Range("E14") "  -0.28%  "

Set-TextValue $ws.Range("D15") "2.234.41"
Set-TextValue $ws.Range("E15") "  -0.60%  "

Set-TextValue $ws.Range("D16") "0.826"
Set-TextValue $ws.Range("E16") "  -2.48%  "

Set-TextValue $ws.Range("D17") "13.60"
Set-TextValue $ws.Range("E17") "  -4.52%  "

Set-TextValue $ws.Range("D18") "43.877.10"
Set-TextValue $ws.Range("E18") "  +0.02%  "

Set-TextValue $ws.Range("D19") "0.0₃0962"
Set-TextValue $ws.Range("E19") "  -1.02%  "

Set-TextValue $ws.Range("D20") "12.12"
Set-TextValue $ws.Range("E20") "  -8.72%  "

Set-TextValue $ws.Range("D21") "6.26"
Set-TextValue $ws.Range("E21") "  -2.24%  "

Set-TextValue $ws.Range("D22") "65.08"
Set-TextValue $ws.Range("E22") "  -0.10%  "

Set-TextValue $ws.Range("D23") "236.55"
Set-TextValue $ws.Range("E23") "  +1.50%  "

Set-TextValue $ws.Range("D24") "2.94"
Set-TextValue $ws.Range("E24") "  -4.62%  "

Set-TextValue $ws.Range("D25") "1.96"
Set-TextValue $ws.Range("E25") "  -4.98%  "

Set-TextValue $ws.Range("E26") "  -0.06%  "

Set-TextValue $ws.Range("D27") "9.96"
Set-TextValue $ws.Range("E27") "  -5.71%  "

Set-TextValue $ws.Range("E28") "  -0.64%  "

Set-TextValue $ws.Range("D29") "37.31"
Set-TextValue $ws.Range("E29") "  -2.61%  "

Set-TextValue $ws.Range("D30") "5.99"
Set-TextValue $ws.Range("E30") "  -0.78%  "

Set-TextValue $ws.Range("D31") "19.91"
Set-TextValue $ws.Range("E31") "  -0.51%  "

Set-TextValue $ws.Range("D32") "152.88"
Set-TextValue $ws.Range("E32") "  -3.31%  "

Set-TextValue $ws.Range("D33") "0.0802"
Set-TextValue $ws.Range("E33") "  -4.50%  "

Set-TextValue $ws.Range("E34") "  +4.66%  "

Set-TextValue $ws.Range("D35") "2.57"
Set-TextValue $ws.Range("E35") "  -3.88%  "

Set-TextValue $ws.Range("B36") "Kaspa"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D36") "0.108"
Set-TextValue $ws.Range("E36") "  -3.64%  "

Set-TextValue $ws.Range("B37") "Stellar"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D37") "0.118"
Set-TextValue $ws.Range("E37") "  +0.55%  "

Set-TextValue $ws.Range("D38") "1.79"
Set-TextValue $ws.Range("E38") "  -6.97%  "

Set-TextValue $ws.Range("D39") "15.16"
Set-TextValue $ws.Range("E39") "  -6.03%  "

Set-TextValue $ws.Range("D40") "3.85"
Set-TextValue $ws.Range("E40") "  -6.91%  "

Set-TextValue $ws.Range("D41") "3.35"
Set-TextValue $ws.Range("E41") "  -8.16%  "

Set-TextValue $ws.Range("E42") "  -3.76%  "

Set-TextValue $ws.Range("E43") "  +0.26%  "

Set-TextValue $ws.Range("D44") "1.725.85"
Set-TextValue $ws.Range("E44") "  -1.88%  "

Set-TextValue $ws.Range("D45") "85.33"
Set-TextValue $ws.Range("E45") "  +5.61%  "

Set-TextValue $ws.Range("E46") "  -3.35%  "

Set-TextValue $ws.Range("D47") "100.06"
Set-TextValue $ws.Range("E47") "  -2.74%  "

Set-TextValue $ws.Range("D48") "4.93"
Set-TextValue $ws.Range("E48") "  -4.07%  "

Set-TextValue $ws.Range("D49") "69.30"
Set-TextValue $ws.Range("E49") "  -6.22%  "

Set-TextValue $ws.Range("D50") "8.07"
Set-TextValue $ws.Range("E50") "  -2.11%  "

Set-TextValue $ws.Range("D51") "54.19"
Set-TextValue $ws.Range("E51") "  -4.88%  "
